# Add a non-blank value in the previously blank row 5 (cell C5 = 30)
# and leave the selection on that cell, matching the author's edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C5").Value = 30
$ws.Range("C5").Select() | Out-Null
